# Increase font sizes throughout the resume document.
# Mapping (Word "Font.Size" is in points; OOXML w:sz is in half-points):
#   16pt (sz=32) -> 18pt (sz=36)   Name header
#    9pt (sz=18) -> 10pt (sz=20)   Contact info / body / bullets / dates
#   12pt (sz=24) -> 13pt (sz=26)   Section headers
#   11pt (sz=22) -> 12pt (sz=24)   Job titles / degree lines
#   10pt (sz=20) -> 11pt (sz=22)   Overview paragraph

$d = $word.ActiveDocument

$sizeMap = @{
    16 = 18
    9  = 10
    12 = 13
    11 = 12
    10 = 11
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pRng = $p.Range

    # Exclude the trailing paragraph-mark character so only the run(s) of
    # actual text get their size changed, not the paragraph mark's own
    # character properties (w:pPr/w:rPr).
    $textEnd = $pRng.End - 1
    if ($textEnd -gt $pRng.Start) {
        $rng = $d.Range($pRng.Start, $textEnd)
        $current = $rng.Font.Size
        if ($sizeMap.ContainsKey($current)) {
            $rng.Font.Size = $sizeMap[$current]
        }
    }
}
